$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# A new leave entry ("SL(1-0-0)" / "SL(2-0-0)") is being recorded, splitting
# the formerly-empty June-2023 row into two rows: the existing row 100 is
# edited in place, and a brand-new row is inserted right after it (pushing
# every subsequent month row down by one).
# ---------------------------------------------------------------------------

# Insert a fresh row above the old row 101 - this shifts old rows 101-148
# down to 102-149, carrying their values/formulas with them.
$ws.Rows.Item(101).Insert()

# The newly-inserted row 101 comes back from Insert() with default (no
# border) formatting. Re-apply the same cell formatting used throughout the
# table by copying it down from row 100.
$ws.Range("A100:K100").Copy()
$ws.Range("A101:K101").PasteSpecial(-4122)

# Grow the table (ListObject) so it covers the new row.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K149"))

# --- Row 100: existing June-2023 row, now carrying the first SL entry -----
$ws.Range("B100").Value = "SL(1-0-0)"
$ws.Range("G100").ClearContents()
$ws.Range("H100").Value = 1
$ws.Range("K100").NumberFormat = "m/d/yy"
$ws.Range("K100").Value = 45077

# --- Row 101: brand-new row carrying the second SL entry ------------------
$ws.Range("B101").Value = "SL(2-0-0)"
$ws.Range("G101").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("H101").Value = 2
$ws.Range("K101").NumberFormat = "m/d/yy"
$ws.Range("K101").Value = "6/20,27/2023"

# The calculated-column formula on the row that used to be 148 (now 149,
# the table's styled "last row") gets mangled by the row-shift; restore the
# original Table1-qualified form so it keeps evaluating to "".
$ws.Range("G149").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Update the selected cell to match where the edit left off.
$ws.Range("K102").Select()
